$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("welk001")
$ws.Range("A1").Value = "test"
